$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed cryptocurrency market data.
# Numeric-looking text values are written via a NumberFormat/ClearFormats
# round-trip so Excel keeps them as text (matching the source workbook)
# instead of silently converting them to floating point numbers.

$ws.Range("D2").Value = "71.398.01"
$ws.Range("E2").Value = "  +7.72%  "
$ws.Range("D3").Value = "3.673.21"
$ws.Range("E3").Value = "  +19.64%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.38%  "
$ws.Range("D7").Value = "3.673.65"
$ws.Range("E7").Value = "  +19.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +5.41%  "
$ws.Range("E10").Value = "  +10.19%  "
$ws.Range("E11").Value = "  +4.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.502"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.89"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +11.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.89%  "
$ws.Range("D15").Value = "4.283.73"
$ws.Range("E15").Value = "  +19.59%  "
$ws.Range("D16").Value = "71.394.22"
$ws.Range("E16").Value = "  +7.79%  "
$ws.Range("D17").Value = "3.672.13"
$ws.Range("E17").Value = "  +19.52%  "
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +8.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "514.44"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +21.25%  "
$ws.Range("E23").Value = "  +10.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.63"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.74%  "
$ws.Range("E26").Value = "  +10.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.87"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.55"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +13.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.33"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +17.07%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000111"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +21.09%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.54%  "
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +11.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.02"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.338"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +12.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "47.69"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.14"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.43%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.09"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.26%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.129"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.92"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +8.57%  "
$ws.Range("D44").Value = "3.156.04"
$ws.Range("E45").Value = "  +11.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "407.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +11.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0367"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +16.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +16.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.56%  "
